$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.306.80'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.626.11'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.94%  '
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '302.39'
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3754'
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '52.48'
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3601'
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.226'
$ws.Range('E10').Value = '  -3.86%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08042'
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.59'
$ws.Range('E13').Value = '  -2.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.530'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001246'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.189'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '1.628.12'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.28'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06912'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.86'
$ws.Range('E20').Value = '  -2.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.436'
$ws.Range('E22').Value = '  -2.15%  '
$ws.Range('D23').Value = '23.308.14'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.65'
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.204'
$ws.Range('E25').Value = '  +2.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.423'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.35'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.296'
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.32'
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.307'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').Value = '1.810.15'
$ws.Range('E32').Value = '  +1.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.752'
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.84'
$ws.Range('E34').Value = '  +3.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9453'
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02820'
$ws.Range('E36').Value = '  +1.12%  '
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.118'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07092'
$ws.Range('E40').Value = '  -4.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.362'
$ws.Range('E41').Value = '  -3.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6994'
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.18'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.26'
$ws.Range('E44').Value = '  -3.09%  '
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.304'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07972'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.201'
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.75'
$ws.Range('E51').Value = '  -4.09%  '
